$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.988.41"
$ws.Range("E2").Value = "  -5.30%  "
$ws.Range("D3").Value = "1.824.55"
$ws.Range("E3").Value = "  -4.11%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").Value = "'328.86"
$ws.Range("E5").Value = "  -2.93%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "'0.4635"
$ws.Range("E7").Value = "  -2.78%  "
$ws.Range("D8").Value = "'0.3843"
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("D9").Value = "'46.00"
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("D11").Value = "'0.9582"
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("D12").Value = "'21.84"
$ws.Range("E12").Value = "  -5.82%  "
$ws.Range("D13").Value = "1.794.51"
$ws.Range("E13").Value = "  -4.64%  "
$ws.Range("D14").Value = "'5.643"
$ws.Range("E14").Value = "  -4.87%  "
$ws.Range("D15").Value = "'6.857"
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("D16").Value = "'0.06852"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "'0.9993"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "'86.35"
$ws.Range("E18").Value = "  -3.15%  "
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("D20").Value = "'16.62"
$ws.Range("E20").Value = "  -4.07%  "
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").Value = "28.007.54"
$ws.Range("E22").Value = "  -5.21%  "
$ws.Range("D23").Value = "'5.299"
$ws.Range("E23").Value = "  -3.66%  "
$ws.Range("D25").Value = "'2.097"
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("D26").Value = "2.038.52"
$ws.Range("E26").Value = "  -3.08%  "
$ws.Range("D27").Value = "'152.44"
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("D28").Value = "'19.15"
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("D29").Value = "'5.739"
$ws.Range("E29").Value = "  -11.62%  "
$ws.Range("D30").Value = "'1.969"
$ws.Range("E30").Value = "  -4.27%  "
$ws.Range("D31").Value = "'116.51"
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("E32").Value = "  -6.04%  "
$ws.Range("D33").Value = "'0.09230"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").Value = "'5.280"
$ws.Range("E34").Value = "  -3.48%  "
$ws.Range("D35").Value = "'1.313"
$ws.Range("E35").Value = "  -5.43%  "
$ws.Range("D36").Value = "'3.338"
$ws.Range("E36").Value = "  -5.57%  "
$ws.Range("D37").Value = "'0.05921"
$ws.Range("E37").Value = "  -8.04%  "
$ws.Range("D39").Value = "'1.140"
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("D40").Value = "'0.9992"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("D42").Value = "'0.5573"
$ws.Range("E42").Value = "  -4.22%  "
$ws.Range("D43").Value = "'9.904"
$ws.Range("E43").Value = "  -6.05%  "
$ws.Range("D44").Value = "'0.1762"
$ws.Range("E44").Value = "  -3.15%  "
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("D46").Value = "'2.196"
$ws.Range("E46").Value = "  -10.38%  "
$ws.Range("D47").Value = "'11.65"
$ws.Range("E47").Value = "  -4.48%  "
$ws.Range("D48").Value = "'0.5246"
$ws.Range("E48").Value = "  -4.15%  "
$ws.Range("D49").Value = "'0.06997"
$ws.Range("E49").Value = "  -5.44%  "
$ws.Range("D50").Value = "'1.820"
$ws.Range("E50").Value = "  -7.05%  "
$ws.Range("D51").Value = "'111.77"
$ws.Range("E51").Value = "  -3.69%  "
